$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds text-formatted dates (e.g. "01-08-2021"). Assigning a
# date-like string directly via .Value causes Excel to auto-convert it
# into a real date serial number, which is not what the source data
# contains (plain text pulled from Banco Central / MV). To keep it as
# literal text (matching the existing rows), build the text via a
# formula in a scratch cell, copy it, and paste-special as values only -
# this preserves the "General"-formatted, shared-string text cell
# without introducing any new number formatting / styles.
$scratch = $ws.Range("Z1")
$scratch.Formula = '="01-09-2021"'
$scratch.Copy()
$ws.Range("A46").PasteSpecial(-4163)
$scratch.Clear()

$ws.Range("B46").Value = 111.45
$ws.Range("C46").Value = 109.43
$ws.Range("D46").Value = 113.37
$ws.Range("E46").Value = 109.18
$ws.Range("F46").Value = 122.49
